# Update test configurations: enable headless mode by flipping the
# ExecutionFlag column on the AddProduct sheet from "No" to "Yes" for
# every existing test case row, and move the active selection.

$wb = $excel.ActiveWorkbook

$wsAdd = $wb.Worksheets.Item("AddProduct")

# Enable every test case (rows 3-8, column A = ExecutionFlag)
for ($r = 3; $r -le 8; $r++) {
    $wsAdd.Cells.Item($r, 1).Value = "Yes"
}

# Move the active selection on the AddProduct sheet
$wsAdd.Range("E14").Select()
